$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.842.92"
$ws.Range("E2").Value = "  +4.10%  "
$ws.Range("D3").Value = "3.072.42"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'578.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.77%  "
$ws.Range("D6").Value = "'141.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.92%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.060.47"
$ws.Range("E8").Value = "  +2.57%  "
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("E10").Value = "  +5.15%  "
$ws.Range("D11").Value = "'5.67"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.95%  "
$ws.Range("E13").Value = "  +3.95%  "
$ws.Range("D14").Value = "'35.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.43%  "
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").Value = "3.581.11"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("D17").Value = "'7.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").Value = "3.069.37"
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("D19").Value = "61.780.70"
$ws.Range("E19").Value = "  +4.18%  "
$ws.Range("D20").Value = "'446.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.88%  "
$ws.Range("D21").Value = "'13.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").Value = "'7.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.24%  "
$ws.Range("D24").Value = "'13.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("D25").Value = "'81.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E27").Value = "  +5.35%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("E29").Value = "  +4.62%  "
$ws.Range("D30").Value = "'8.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.44%  "
$ws.Range("D31").Value = "'6.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.61%  "
$ws.Range("D32").Value = "'0.112"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +13.23%  "
$ws.Range("D33").Value = "'26.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.12%  "
$ws.Range("D34").Value = "'1.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.08%  "
$ws.Range("D35").Value = "0.0₃0788"
$ws.Range("E35").Value = "  +2.93%  "
$ws.Range("D36").Value = "'6.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.80%  "
$ws.Range("E37").Value = "  +4.88%  "
$ws.Range("D38").Value = "'49.97"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("D39").Value = "'2.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.70%  "
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("D41").Value = "'420.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.57%  "
$ws.Range("D42").Value = "2.978.03"
$ws.Range("E42").Value = "  +8.12%  "
$ws.Range("D43").Value = "'0.0369"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.91%  "
$ws.Range("E44").Value = "  +9.77%  "
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("D46").Value = "'2.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.76%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'124.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.08%  "
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").Value = "'34.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "'24.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.51%  "

Write-Output "Applied 84 cell updates"
